$wb = $excel.ActiveWorkbook

# --- Sheet1 (Step1_Data): update raw signal-value cells ---
$wsData = $wb.Worksheets.Item("Step1_Data")
$wsData.Range("D2").Value = 0.04474526330661761
$wsData.Range("E2").Value = 0.04679875640740945
$wsData.Range("F2").Value = 0.1858191962305922
$wsData.Range("G2").Value = 0.1125177710422336
$wsData.Range("J2").Value = 0.009307046095662573
$wsData.Range("L2").Value = 0.02302389220976572
$wsData.Range("M2").Value = 0.04649264580020804
$wsData.Range("N2").Value = 0.02159071061666533
$wsData.Range("O2").Value = 0.2041287989386799
$wsData.Range("R2").Value = 0.03619464797014928
$wsData.Range("S2").Value = 0.04597225083647536
$wsData.Range("T2").Value = 0.09749425458540892
$wsData.Range("U2").Value = 0.01819897956677105
$wsData.Range("V2").Value = 0.004388487241468876
$wsData.Range("W2").Value = 0.002808437144182921
$wsData.Range("AA2").Value = 0.005604237759174301
$wsData.Range("AD2").Value = 0.03452656172554074
$wsData.Range("AF2").Value = 0.0507951089146542
$wsData.Range("AG2").Value = 0.00954096345412989
$wsData.Range("AH2").Value = 0.00005199015421018281
$wsData.Range("D3").Value = 0.08884610831446711
$wsData.Range("E3").Value = 0.02538123859597912
$wsData.Range("F3").Value = 0.1627697398449063
$wsData.Range("G3").Value = 0.09274924673759515
$wsData.Range("J3").Value = 0.03028478303538525
$wsData.Range("L3").Value = 0.02871610530597796
$wsData.Range("M3").Value = 0.04310191185983549
$wsData.Range("N3").Value = 0.01437565964664485
$wsData.Range("O3").Value = 0.2539131776173248
$wsData.Range("R3").Value = 0.02269224746840743
$wsData.Range("S3").Value = 0.0397225361630474
$wsData.Range("T3").Value = 0.09647575750643887
$wsData.Range("V3").Value = 0.001313778459563746
$wsData.Range("AA3").Value = 0.02702538586481252
$wsData.Range("AD3").Value = 0.007021641827594435
$wsData.Range("AF3").Value = 0.0339781168321098
$wsData.Range("AG3").Value = 0.03163256491990968
$wsData.Range("D4").Value = 0.2616502453699902
$wsData.Range("E4").Value = 0.006604330012431536
$wsData.Range("F4").Value = 0.2874220780238156
$wsData.Range("G4").Value = 0.01308147622289825
$wsData.Range("H4").Value = 0.00441175496244287
$wsData.Range("I4").Value = 0.006930072086128379
$wsData.Range("L4").Value = 0.05149651867366712
$wsData.Range("N4").Value = 0.1279079413931333
$wsData.Range("O4").Value = 0.008263821974867849
$wsData.Range("Q4").Value = 0.02229611864536993
$wsData.Range("S4").Value = 0.1025050263053454
$wsData.Range("U4").Value = 0.0217903123731341
$wsData.Range("Z4").Value = 0.01733417283156081
$wsData.Range("AC4").Value = 0.004350476999973668
$wsData.Range("AD4").Value = 0.01860537790183923
$wsData.Range("AF4").Value = 0.04535027622340158
$wsData.Range("D5").Value = 0.226287662041849
$wsData.Range("F5").Value = 0.1748442556992893
$wsData.Range("G5").Value = 0.02632773609184789
$wsData.Range("J5").Value = 0.02316078329112299
$wsData.Range("L5").Value = 0.07571771768274636
$wsData.Range("M5").Value = 0.005393393348952774
$wsData.Range("N5").Value = 0.1126443933341342
$wsData.Range("O5").Value = 0.09407710142400329
$wsData.Range("Q5").Value = 0.00421551244140142
$wsData.Range("S5").Value = 0.103172248246469
$wsData.Range("T5").Value = 0.03207822532985424
$wsData.Range("U5").Value = 0.02749791974501179
$wsData.Range("W5").Value = 0.004611587075943744
$wsData.Range("Z5").Value = 0.007446146702510453
$wsData.Range("AA5").Value = 0.006028214550166459
$wsData.Range("AD5").Value = 0.01038565235006716
$wsData.Range("AF5").Value = 0.06611145064463025
$wsData.Range("D6").Value = 0.005718516937492077
$wsData.Range("E6").Value = 0.1260239386286452
$wsData.Range("F6").Value = 0.09703236053234364
$wsData.Range("G6").Value = 0.2516959581816063
$wsData.Range("I6").Value = 0.002870632600202734
$wsData.Range("L6").Value = 0.003723239759171308
$wsData.Range("M6").Value = 0.03697724671061999
$wsData.Range("O6").Value = 0.1984003352004571
$wsData.Range("R6").Value = 0.03550343160022515
$wsData.Range("T6").Value = 0.1232181012457643
$wsData.Range("V6").Value = 0.01316348150469672
$wsData.Range("AA6").Value = 0.01766478596091157
$wsData.Range("AC6").Value = 0.005941462957936017
$wsData.Range("AD6").Value = 0.01415939195377156
$wsData.Range("AE6").Value = 0.007975973097530195
$wsData.Range("AF6").Value = 0.02082942988389696
$wsData.Range("AG6").Value = 0.03910171324472909

# --- Sheet2 (Step2_Sj): running cumulative sum per row, recomputed from Step1_Data ---
$wsSj = $wb.Worksheets.Item("Step2_Sj")
$wsSj.Range("B2").Value = 0
$wsSj.Range("C2").Value = 0
$wsSj.Range("D2").Value = 0.04474526330661761
$wsSj.Range("E2").Value = 0.09154401971402706
$wsSj.Range("F2").Value = 0.27736321594461927
$wsSj.Range("G2").Value = 0.38988098698685286
$wsSj.Range("H2").Value = 0.38988098698685286
$wsSj.Range("I2").Value = 0.38988098698685286
$wsSj.Range("J2").Value = 0.39918803308251544
$wsSj.Range("K2").Value = 0.39918803308251544
$wsSj.Range("L2").Value = 0.42221192529228113
$wsSj.Range("M2").Value = 0.4687045710924892
$wsSj.Range("N2").Value = 0.49029528170915454
$wsSj.Range("O2").Value = 0.6944240806478345
$wsSj.Range("P2").Value = 0.6944240806478345
$wsSj.Range("Q2").Value = 0.6944240806478345
$wsSj.Range("R2").Value = 0.7306187286179837
$wsSj.Range("S2").Value = 0.7765909794544591
$wsSj.Range("T2").Value = 0.874085234039868
$wsSj.Range("U2").Value = 0.8922842136066391
$wsSj.Range("V2").Value = 0.896672700848108
$wsSj.Range("W2").Value = 0.8994811379922909
$wsSj.Range("X2").Value = 0.8994811379922909
$wsSj.Range("Y2").Value = 0.8994811379922909
$wsSj.Range("Z2").Value = 0.8994811379922909
$wsSj.Range("AA2").Value = 0.9050853757514652
$wsSj.Range("AB2").Value = 0.9050853757514652
$wsSj.Range("AC2").Value = 0.9050853757514652
$wsSj.Range("AD2").Value = 0.9396119374770059
$wsSj.Range("AE2").Value = 0.9396119374770059
$wsSj.Range("AF2").Value = 0.9904070463916601
$wsSj.Range("AG2").Value = 0.99994800984579
$wsSj.Range("AH2").Value = 1.0000000000000002
$wsSj.Range("AI2").Value = 1.0000000000000002
$wsSj.Range("B3").Value = 0
$wsSj.Range("C3").Value = 0
$wsSj.Range("D3").Value = 0.08884610831446711
$wsSj.Range("E3").Value = 0.11422734691044623
$wsSj.Range("F3").Value = 0.2769970867553525
$wsSj.Range("G3").Value = 0.36974633349294767
$wsSj.Range("H3").Value = 0.36974633349294767
$wsSj.Range("I3").Value = 0.36974633349294767
$wsSj.Range("J3").Value = 0.4000311165283329
$wsSj.Range("K3").Value = 0.4000311165283329
$wsSj.Range("L3").Value = 0.42874722183431085
$wsSj.Range("M3").Value = 0.4718491336941463
$wsSj.Range("N3").Value = 0.48622479334079116
$wsSj.Range("O3").Value = 0.740137970958116
$wsSj.Range("P3").Value = 0.740137970958116
$wsSj.Range("Q3").Value = 0.740137970958116
$wsSj.Range("R3").Value = 0.7628302184265234
$wsSj.Range("S3").Value = 0.8025527545895709
$wsSj.Range("T3").Value = 0.8990285120960098
$wsSj.Range("U3").Value = 0.8990285120960098
$wsSj.Range("V3").Value = 0.9003422905555735
$wsSj.Range("W3").Value = 0.9003422905555735
$wsSj.Range("X3").Value = 0.9003422905555735
$wsSj.Range("Y3").Value = 0.9003422905555735
$wsSj.Range("Z3").Value = 0.9003422905555735
$wsSj.Range("AA3").Value = 0.927367676420386
$wsSj.Range("AB3").Value = 0.927367676420386
$wsSj.Range("AC3").Value = 0.927367676420386
$wsSj.Range("AD3").Value = 0.9343893182479804
$wsSj.Range("AE3").Value = 0.9343893182479804
$wsSj.Range("AF3").Value = 0.9683674350800903
$wsSj.Range("AG3").Value = 1
$wsSj.Range("AH3").Value = 1
$wsSj.Range("AI3").Value = 1
$wsSj.Range("B4").Value = 0
$wsSj.Range("C4").Value = 0
$wsSj.Range("D4").Value = 0.2616502453699902
$wsSj.Range("E4").Value = 0.26825457538242176
$wsSj.Range("F4").Value = 0.5556766534062374
$wsSj.Range("G4").Value = 0.5687581296291356
$wsSj.Range("H4").Value = 0.5731698845915785
$wsSj.Range("I4").Value = 0.5800999566777069
$wsSj.Range("J4").Value = 0.5800999566777069
$wsSj.Range("K4").Value = 0.5800999566777069
$wsSj.Range("L4").Value = 0.631596475351374
$wsSj.Range("M4").Value = 0.631596475351374
$wsSj.Range("N4").Value = 0.7595044167445073
$wsSj.Range("O4").Value = 0.7677682387193752
$wsSj.Range("P4").Value = 0.7677682387193752
$wsSj.Range("Q4").Value = 0.7900643573647451
$wsSj.Range("R4").Value = 0.7900643573647451
$wsSj.Range("S4").Value = 0.8925693836700905
$wsSj.Range("T4").Value = 0.8925693836700905
$wsSj.Range("U4").Value = 0.9143596960432246
$wsSj.Range("V4").Value = 0.9143596960432246
$wsSj.Range("W4").Value = 0.9143596960432246
$wsSj.Range("X4").Value = 0.9143596960432246
$wsSj.Range("Y4").Value = 0.9143596960432246
$wsSj.Range("Z4").Value = 0.9316938688747854
$wsSj.Range("AA4").Value = 0.9316938688747854
$wsSj.Range("AB4").Value = 0.9316938688747854
$wsSj.Range("AC4").Value = 0.9360443458747592
$wsSj.Range("AD4").Value = 0.9546497237765984
$wsSj.Range("AE4").Value = 0.9546497237765984
$wsSj.Range("AF4").Value = 1
$wsSj.Range("AG4").Value = 1
$wsSj.Range("AH4").Value = 1
$wsSj.Range("AI4").Value = 1
$wsSj.Range("B5").Value = 0
$wsSj.Range("C5").Value = 0
$wsSj.Range("D5").Value = 0.226287662041849
$wsSj.Range("E5").Value = 0.226287662041849
$wsSj.Range("F5").Value = 0.4011319177411383
$wsSj.Range("G5").Value = 0.42745965383298623
$wsSj.Range("H5").Value = 0.42745965383298623
$wsSj.Range("I5").Value = 0.42745965383298623
$wsSj.Range("J5").Value = 0.4506204371241092
$wsSj.Range("K5").Value = 0.4506204371241092
$wsSj.Range("L5").Value = 0.5263381548068555
$wsSj.Range("M5").Value = 0.5317315481558083
$wsSj.Range("N5").Value = 0.6443759414899425
$wsSj.Range("O5").Value = 0.7384530429139458
$wsSj.Range("P5").Value = 0.7384530429139458
$wsSj.Range("Q5").Value = 0.7426685553553471
$wsSj.Range("R5").Value = 0.7426685553553471
$wsSj.Range("S5").Value = 0.8458408036018161
$wsSj.Range("T5").Value = 0.8779190289316704
$wsSj.Range("U5").Value = 0.9054169486766822
$wsSj.Range("V5").Value = 0.9054169486766822
$wsSj.Range("W5").Value = 0.910028535752626
$wsSj.Range("X5").Value = 0.910028535752626
$wsSj.Range("Y5").Value = 0.910028535752626
$wsSj.Range("Z5").Value = 0.9174746824551364
$wsSj.Range("AA5").Value = 0.9235028970053029
$wsSj.Range("AB5").Value = 0.9235028970053029
$wsSj.Range("AC5").Value = 0.9235028970053029
$wsSj.Range("AD5").Value = 0.93388854935537
$wsSj.Range("AE5").Value = 0.93388854935537
$wsSj.Range("AF5").Value = 1.0000000000000002
$wsSj.Range("AG5").Value = 1.0000000000000002
$wsSj.Range("AH5").Value = 1.0000000000000002
$wsSj.Range("AI5").Value = 1.0000000000000002
$wsSj.Range("B6").Value = 0
$wsSj.Range("C6").Value = 0
$wsSj.Range("D6").Value = 0.005718516937492077
$wsSj.Range("E6").Value = 0.13174245556613728
$wsSj.Range("F6").Value = 0.22877481609848094
$wsSj.Range("G6").Value = 0.4804707742800872
$wsSj.Range("H6").Value = 0.4804707742800872
$wsSj.Range("I6").Value = 0.48334140688029
$wsSj.Range("J6").Value = 0.48334140688029
$wsSj.Range("K6").Value = 0.48334140688029
$wsSj.Range("L6").Value = 0.4870646466394613
$wsSj.Range("M6").Value = 0.5240418933500813
$wsSj.Range("N6").Value = 0.5240418933500813
$wsSj.Range("O6").Value = 0.7224422285505384
$wsSj.Range("P6").Value = 0.7224422285505384
$wsSj.Range("Q6").Value = 0.7224422285505384
$wsSj.Range("R6").Value = 0.7579456601507636
$wsSj.Range("S6").Value = 0.7579456601507636
$wsSj.Range("T6").Value = 0.8811637613965279
$wsSj.Range("U6").Value = 0.8811637613965279
$wsSj.Range("V6").Value = 0.8943272429012246
$wsSj.Range("W6").Value = 0.8943272429012246
$wsSj.Range("X6").Value = 0.8943272429012246
$wsSj.Range("Y6").Value = 0.8943272429012246
$wsSj.Range("Z6").Value = 0.8943272429012246
$wsSj.Range("AA6").Value = 0.9119920288621362
$wsSj.Range("AB6").Value = 0.9119920288621362
$wsSj.Range("AC6").Value = 0.9179334918200721
$wsSj.Range("AD6").Value = 0.9320928837738437
$wsSj.Range("AE6").Value = 0.9400688568713739
$wsSj.Range("AF6").Value = 0.9608982867552708
$wsSj.Range("AG6").Value = 0.9999999999999999
$wsSj.Range("AH6").Value = 0.9999999999999999
$wsSj.Range("AI6").Value = 0.9999999999999999

# --- Sheet3-6 (Step3_DataPts_*): recompute threshold crossing point (D,F,G) ---
$ws3 = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws3.Range("D2").Value = 14
$ws3.Range("F2").Value = 0.6944240806478345
$ws3.Range("G2").Value = 12
$ws3.Range("D3").Value = 14
$ws3.Range("F3").Value = 0.740137970958116
$ws3.Range("G3").Value = 13
$ws3.Range("D4").Value = 5
$ws3.Range("F4").Value = 0.5556766534062374
$ws3.Range("G4").Value = 4
$ws3.Range("D5").Value = 11
$ws3.Range("F5").Value = 0.5263381548068555
$ws3.Range("G5").Value = 10
$ws3.Range("D6").Value = 12
$ws3.Range("F6").Value = 0.5240418933500813
$ws3.Range("G6").Value = 10
$ws4 = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws4.Range("D2").Value = 17
$ws4.Range("F2").Value = 0.7306187286179837
$ws4.Range("G2").Value = 15
$ws4.Range("D3").Value = 14
$ws4.Range("F3").Value = 0.740137970958116
$ws4.Range("G3").Value = 13
$ws4.Range("D4").Value = 13
$ws4.Range("F4").Value = 0.7595044167445073
$ws4.Range("G4").Value = 12
$ws4.Range("D5").Value = 14
$ws4.Range("F5").Value = 0.7384530429139458
$ws4.Range("G5").Value = 13
$ws4.Range("D6").Value = 14
$ws4.Range("F6").Value = 0.7224422285505384
$ws4.Range("G6").Value = 12
$ws5 = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws5.Range("D2").Value = 19
$ws5.Range("F2").Value = 0.874085234039868
$ws5.Range("G2").Value = 17
$ws5.Range("D3").Value = 18
$ws5.Range("F3").Value = 0.8025527545895709
$ws5.Range("G3").Value = 17
$ws5.Range("D4").Value = 18
$ws5.Range("F4").Value = 0.8925693836700905
$ws5.Range("G4").Value = 17
$ws5.Range("D5").Value = 18
$ws5.Range("F5").Value = 0.8458408036018161
$ws5.Range("G5").Value = 17
$ws5.Range("D6").Value = 19
$ws5.Range("F6").Value = 0.8811637613965279
$ws5.Range("G6").Value = 17
$ws6 = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws6.Range("D2").Value = 26
$ws6.Range("F2").Value = 0.9050853757514652
$ws6.Range("G2").Value = 24
$ws6.Range("D3").Value = 21
$ws6.Range("F3").Value = 0.9003422905555735
$ws6.Range("G3").Value = 20
$ws6.Range("D4").Value = 20
$ws6.Range("F4").Value = 0.9143596960432246
$ws6.Range("G4").Value = 19
$ws6.Range("D5").Value = 20
$ws6.Range("F5").Value = 0.9054169486766822
$ws6.Range("G5").Value = 19
$ws6.Range("D6").Value = 26
$ws6.Range("F6").Value = 0.9119920288621362
$ws6.Range("G6").Value = 24

# --- Tire_Type label correction: "710Rバフ50" -> "710R" (K column) ---
$ws3.Range("K2:K6").Value = "710R"
$ws4.Range("K2:K6").Value = "710R"
$ws5.Range("K2:K6").Value = "710R"
$ws6.Range("K2:K6").Value = "710R"
